$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Typo fixes: "Which" -> "Whitch" in the Q6 / Q7 questions only
#    (other "Which" occurrences elsewhere in the doc must stay as-is)
# ---------------------------------------------------------------
$d.Content.Find.Execute("Q6 Which are faster", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Q6 Whitch are faster", 2) | Out-Null

$d.Content.Find.Execute("Q7 Which are faster", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Q7 Whitch are faster", 2) | Out-Null

# ---------------------------------------------------------------
# 2) Replace the "queries_3.ipynb" block's old Q2 question with four
#    new questions (Q2-Q5) followed by two blank paragraphs, and
#    blank out the old paragraph that used to hold the Q2 text
#    (its paragraph mark / formatting is preserved, just like the
#    source revision keeps an empty trailing paragraph there).
# ---------------------------------------------------------------
$oldQ2Text = "Q2 Can be (the best most field GOAL) guarantee to be (the best winner)?`r"

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq $oldQ2Text) {
        $target = $p
    }
}

if ($target -ne $null) {
    $insertRange = $target.Range
    $newText = "Q2 Select top 10 teams with its wining percentage and with total points?`r" + `
               "Q3 In this Query group data by wining percentage?`r" + `
               "Q4 In this Query group data by Team Name?`r" + `
               "Q5 In this Query Reduce by key and calculate the average by the key?`r" + `
               "`r" + `
               "`r"
    $insertRange.InsertBefore($newText)

    # InsertBefore() re-seats the original $target paragraph object onto
    # the first newly-inserted paragraph, so re-locate the paragraph that
    # still holds the old sentence (now six paragraphs further down) and
    # clear its text while keeping the paragraph mark/formatting intact.
    $oldParagraph = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq $oldQ2Text) {
            $oldParagraph = $p
        }
    }

    if ($oldParagraph -ne $null) {
        $oldRange = $oldParagraph.Range
        $clearRange = $d.Range($oldRange.Start, $oldRange.End - 1)
        $clearRange.Text = ""
    }
}

Write-Output "done"
